$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the checklist answers that were entered in the wrong rows:
# B36 ("NA") and B37 (the long Edinburgh Postnatal Depression Scale note)
# need to trade places.
$b36 = $ws.Range("B36").Value2
$b37 = $ws.Range("B37").Value2

$ws.Range("B36").Value = $b37
$ws.Range("B37").Value = $b36

# B37 had been highlighted red from a previous (incorrect) entry; now that
# it just holds "NA" again, clear the highlight back to the normal
# (theme background / white) fill.
$ws.Range("B37").Interior.ThemeColor = 2
$ws.Range("B37").Interior.TintAndShade = 0
